# Generate Report for Handback
#
# Refreshes the "Latest HO Xliff Generate Date" / handoff / handback
# timestamps that are stamped whenever the handback report is (re-)generated
# for the file "6fcbf077-f894-4b9d-a26f-a248bbfe555d.md" (row 2 of each
# per-locale sheet). The companion file "79e2398d-...md" (row 3) was not
# regenerated in this run, so its timestamps are left untouched.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 02:54:20"

# --- zh-cn handback detail sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-05 02:54:15"
$wsZhCn.Range("K2").Value = "2016-09-05 02:54:31"

# --- de-de handback detail sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-05 02:54:20"
$wsDeDe.Range("K2").Value = "2016-09-05 02:54:39"
